# "Separate problems from solutions"
#
# The deck originally held 4 slides: a "Problem" divider, the "Sudoku"
# slide (with the puzzle/solution pictures), a "Solution" divider, and a
# content slide with three more pictures + a caption. This edit keeps only
# the "Sudoku" slide, retitles it "Sudoku Problems", and drops the other
# three slides. It also refreshes the cached datetimeFigureOut placeholder
# text (master + every layout) from 3/25/2023 to 3/30/2023.

$p = $ppt.ActivePresentation

# --- Drop the "Problem" divider, "Solution" divider, and the
#     picture/caption slide, keeping only the "Sudoku" slide. ---
# Work from the back of the deck forward so earlier indices stay stable.
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $s = $p.Slides.Item($i)
    $keep = $false
    for ($j = 1; $j -le $s.Shapes.Count; $j++) {
        $sh = $s.Shapes.Item($j)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            if ($sh.TextFrame.TextRange.Text -eq "Sudoku") {
                $keep = $true
            }
        }
    }
    if (-not $keep) {
        $s.Delete()
    }
}

# --- Retitle the remaining slide. ---
$s = $p.Slides.Item(1)
for ($j = 1; $j -le $s.Shapes.Count; $j++) {
    $sh = $s.Shapes.Item($j)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.TextFrame.TextRange.Text -eq "Sudoku") {
            $sh.TextFrame.TextRange.Text = "Sudoku Problems"
        }
    }
}

# --- Refresh the cached "datetimeFigureOut" field text: 3/25/2023 -> 3/30/2023.
#     These live on the slide master and on every slide layout's Date
#     placeholder (the placeholder shape's ordinal/name varies per layout,
#     so match by name prefix rather than a fixed shape index). Go through
#     $p.SlideMaster (not Slide.Master / Master.CustomLayouts.Item) since
#     that is what reliably addresses each of the 17 distinct layouts. ---
function Update-CachedDate($shapes) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shape = $shapes.Item($k)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            if ($shape.Name -like "Date Placeholder*") {
                if ($shape.TextFrame.TextRange.Text -eq "3/25/2023") {
                    $shape.TextFrame.TextRange.Text = "3/30/2023"
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-CachedDate $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-CachedDate $layout.Shapes
}
